# Modify Xpert for MDR identification
# Insert a new parameter row "int_prop_xpert_sensitivity_mdr" = 0.94399999999999995
# above the existing "program_prop_treatment_success_mdr_asds" row (row 61) on the
# "constants" sheet, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row at row 61 (this shifts rows 61:260 down to 62:261 and copies
# formatting from the row above, matching styles s=26/27/27/27/24 on A:E).
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row.
$ws.Range("A61").Value = "int_prop_xpert_sensitivity_mdr"
$ws.Range("B61").Value = 0.94399999999999995

# Restore the selection/view state as closely as possible: the author's last
# selection after the edit was the full new row 61 (A61:XFD61).
$ws.Range("A61:XFD61").Select()
